# Updates cryptos list values (Price / Volume(1h) columns) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values are plain decimal numbers (e.g. "591.22"). Excel auto-converts such
# strings to numeric values on assignment, which would lose the exact textual formatting
# (trailing zeros, etc.) used throughout this sheet. Force those specific cells to Text
# format first so the values are stored exactly as given.
$textForceCells = @("D5", "D11", "D16", "D18", "D19", "D20", "D21", "D22", "D24", "D36", "D37", "D38", "D40", "D45", "D48", "D49")
foreach ($ref in $textForceCells) {
  $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.372.88"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.507.17"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "591.22"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("E9").Value = "  +5.78%  "
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D11").Value = "0.390"
$ws.Range("E11").Value = "  +3.68%  "
$ws.Range("D12").Value = "4.104.16"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "3.506.13"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "25.80"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "64.358.86"
$ws.Range("D18").Value = "10.07"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").Value = "5.76"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").Value = "13.63"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").Value = "393.30"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").Value = "0.582"
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("D23").Value = "3.645.87"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").Value = "74.46"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -1.02%  "
$ws.Range("E27").Value = "  +3.52%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("E32").Value = "  -5.44%  "
$ws.Range("E33").Value = "  +7.53%  "
$ws.Range("D34").Value = "3.532.74"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "23.39"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "5.35"
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").Value = "6.97"
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("E39").Value = "  +2.21%  "
$ws.Range("D40").Value = "166.26"
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("E41").Value = "  +1.31%  "
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").Value = "24.96"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").Value = "6.81"
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("D49").Value = "0.914"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("D50").Value = "2.380.62"
$ws.Range("E50").Value = "  -3.62%  "
$ws.Range("E51").Value = "  +0.43%  "
